$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25: top-up POST - add estimate values in F and G
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = "12 + "

# Row 26: top-up GET (all) - add estimate values in F and G
$ws.Range("F26").Value = 5
$ws.Range("G26").Value = 7

# Row 27: top-up GET (one) - add estimate values in F and G
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 3

# Update the view state: scrolled down with new selection
$ws.Application.ActiveWindow.ScrollRow = 18
$ws.Range("G32").Select()
